$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.195.97'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.99%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.338.18'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.44%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '582.85'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.57%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '177.10'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +1.94%  '

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.08%  '

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.67%  '

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +3.61%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.582'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +1.33%  '

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +5.91%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000273'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +1.54%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '693.74'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +4.60%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.883.54'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.49%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.42'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.54%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '68.220.92'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.81%  '

$ws.Range("B17").Value = 'TRON'
$ws.Range("C17").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.119'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.35%  '

$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.339.68'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.56%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.43'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.14%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.15'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +2.15%  '

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.91%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.46'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +1.36%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '16.97'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.09%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '100.02'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +1.50%  '

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +2.01%  '

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.95%  '

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +3.05%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '32.95'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -2.09%  '

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +1.06%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.93'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -5.88%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '563.58'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -4.84%  '

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +1.30%  '

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +1.22%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '57.62'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +1.48%  '

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.06%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.691.17'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.08%  '

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +1.14%  '

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +3.81%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '34.77'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +4.31%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.16'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +2.11%  '

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.40%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0₃0670'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +1.44%  '

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.79%  '

$ws.Range("B44").Value = 'ApeXProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.25'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.14%  '

$ws.Range("B45").Value = 'VeChain'
$ws.Range("C45").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0413'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +1.70%  '

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +2.54%  '

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.75%  '

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.17%  '

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.45%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '130.92'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +3.40%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.57'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.76%  '
